$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add a new row to the "Tableau1" table for the backup server
$newRow = $lo.ListRows.Add()

$ws.Range("A7").Value = "srvbackup"
$ws.Range("B7").Value = "VM-Ubunut-serv"
$ws.Range("C7").Value = "192.168.99.110"
$ws.Range("D7").Value = "192.168.99.0/24"

# Rename the "LAN" column header to "VLAN"
$ws.Range("E1").Value = "VLAN"

$ws.Range("E7").Value = 99
$ws.Range("F7").Value = "Backup"

# Copy formatting from the row above onto the new row
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$ws.Range("F6").Select()
